$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.791.72"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "2.119.93"
$ws.Range("E3").Value = "  +10.28%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'255.76"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").Value = "'0.672"
$ws.Range("E6").Value = "  -3.25%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "'48.00"
$ws.Range("E8").Value = "  +9.43%  "
$ws.Range("E9").Value = "  +5.14%  "
$ws.Range("E10").Value = "  +2.42%  "
$ws.Range("E11").Value = "  -3.15%  "
$ws.Range("D13").Value = "2.423.16"
$ws.Range("E13").Value = "  +10.29%  "
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("E15").Value = "  +5.64%  "
$ws.Range("D16").Value = "2.118.11"
$ws.Range("E16").Value = "  +10.46%  "
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").Value = "36.878.85"
$ws.Range("E18").Value = "  +1.13%  "
$ws.Range("D19").Value = "'73.87"
$ws.Range("E19").Value = "  -0.77%  "
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").Value = "'13.47"
$ws.Range("E21").Value = "  +1.92%  "
$ws.Range("D22").Value = "'241.83"
$ws.Range("E22").Value = "  -4.19%  "
$ws.Range("D23").Value = "'5.19"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("E24").Value = "  +0.03%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  -8.80%  "
$ws.Range("D26").Value = "'171.82"
$ws.Range("E26").Value = "  +2.44%  "
$ws.Range("D27").Value = "'21.52"
$ws.Range("E27").Value = "  +14.08%  "
$ws.Range("E28").Value = "  +5.35%  "
$ws.Range("D29").Value = "'2.04"
$ws.Range("E29").Value = "  -7.71%  "
$ws.Range("E30").Value = "  -4.55%  "
$ws.Range("D31").Value = "'25.84"
$ws.Range("E31").Value = "  +55.76%  "
$ws.Range("D32").Value = "'4.51"
$ws.Range("E32").Value = "  -0.43%  "
$ws.Range("D33").Value = "'0.0955"
$ws.Range("E33").Value = "  +11.86%  "
$ws.Range("D34").Value = "'0.0600"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "'2.37"
$ws.Range("E35").Value = "  +18.04%  "
$ws.Range("E36").Value = "  +0.06%  "
$ws.Range("B37").Value = "WEMIXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D37").Value = "'1.88"
$ws.Range("E37").Value = "  -3.23%  "
$ws.Range("B38").Value = "InternetComputer(DFINITY)"
$ws.Range("C38").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D38").Value = "'4.19"
$ws.Range("E38").Value = "  -3.51%  "
$ws.Range("E39").Value = "  +6.86%  "
$ws.Range("E40").Value = "  -7.56%  "
$ws.Range("E41").Value = "  +6.77%  "
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "'99.15"
$ws.Range("E43").Value = "  -7.30%  "
$ws.Range("E44").Value = "  +16.02%  "
$ws.Range("D45").Value = "'16.23"
$ws.Range("E45").Value = "  -4.76%  "
$ws.Range("D46").Value = "1.359.85"
$ws.Range("E46").Value = "  +1.18%  "
$ws.Range("D47").Value = "'0.0841"
$ws.Range("E47").Value = "  +4.16%  "
$ws.Range("D48").Value = "'7.11"
$ws.Range("E48").Value = "  +10.45%  "
$ws.Range("D49").Value = "2.310.05"
$ws.Range("E49").Value = "  +10.20%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  -3.38%  "
$ws.Range("D51").Value = "'2.83"
$ws.Range("E51").Value = "  +1.34%  "
